# Secondary Spacecraft Risk Analysis -- COM risk analysis done
$d = $word.ActiveDocument

# --- Change 1: STR section "Minimal: component selection ..." paragraph --
# Merge the three runs ("...previous" / bookmark / " mission experience")
# into a single run and drop the _GoBack bookmark, leaving the preceding
# "Minimal" / ": " runs untouched.
$rng = $d.Content
$rng.Find.Execute("component selection is limited and some components will come from previous mission experience", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$mergeStart = $rng.Start
$rng.Delete()
$ins = $d.Range($mergeStart, $mergeStart)
$ins.InsertAfter("component selection is limited and some components will come from previous mission experience")

# --- Change 2a: insert new "Minimal: assumptions ..." paragraph right ---
# --- after "Wrong assumptions in link analysis." (no numbering)       ---
$rng = $d.Content
$rng.Find.Execute("Wrong assumptions in link analysis.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint = $d.Range($rng.End, $rng.End)
$xmlMinimal = '<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t>Minimal: assumptions that have been made so far do not have a significant impact on the overall link budget</w:t></w:r></w:p><w:sectPr><w:pgSz w:w="12240" w:h="15840"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($xmlMinimal)

# --- Change 2b: insert new "Low: spending too much time ..." paragraph --
# --- right after "Preliminary component research takes too long"      --
# --- (no numbering, carries the relocated _GoBack bookmark)           --
$rng = $d.Content
$rng.Find.Execute("Preliminary component research takes too long", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint = $d.Range($rng.End, $rng.End)
$xmlLow = '<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t>Low</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>spending too much time on research delaying communication subsystem development</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:sectPr><w:pgSz w:w="12240" w:h="15840"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($xmlLow)

# --- Change 3: "Mistakes in detumbling analysis" -- drop proofErr wrap --
$rng = $d.Content
$rng.Find.Execute("Mistakes in detumbling analysis", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$mergeStart = $rng.Start
$rng.Delete()
$ins = $d.Range($mergeStart, $mergeStart)
$ins.InsertAfter("Mistakes in detumbling analysis")

# --- Change 4: "Underestimated detumbling requirements" -- same cleanup --
$rng = $d.Content
$rng.Find.Execute("Underestimated detumbling requirements", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$mergeStart = $rng.Start
$rng.Delete()
$ins = $d.Range($mergeStart, $mergeStart)
$ins.InsertAfter("Underestimated detumbling requirements")
